$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7:E11").Value = "Completed"

$ws.Range("H11").Select()
